$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignIn")
$ws.Columns.Item(10).ColumnWidth = 22.5
Write-Output "done"
